{"js": "// Commit: \"Update master to output generated at 9a8706d\"\n// The document is a single table of \"A\u00d7B=C\" multiplication-fact cells.\n// This edit swaps 25 of those problem/answer strings for new ones while\n// leaving every other paragraph, row, and run formatting untouched.\n//\n// Each [oldText, newText] pair below is a unique, exact string that occurs\n// exactly once in the document body, so a case-sensitive search-and-replace\n// on the run text (which preserves the run's rFonts/sz formatting) is\n// sufficient and unambiguous.\nconst replacements = [\n  [\"651\u00d77=4557\", \"316\u00d74=1264\"],\n  [\"624\u00d77=4368\", \"190\u00d77=1330\"],\n  [\"144\u00d73=432\", \"152\u00d75=760\"],\n  [\"291\u00d78=2328\", \"129\u00d76=774\"],\n  [\"675\u00d76=4050\", \"212\u00d72=424\"],\n  [\"321\u00d76=1926\", \"959\u00d78=7672\"],\n  [\"200\u00d74=800\", \"514\u00d79=4626\"],\n  [\"389\u00d74=1556\", \"261\u00d74=1044\"],\n  [\"302\u00d78=2416\", \"521\u00d74=2084\"],\n  [\"873\u00d73=2619\", \"120\u00d73=360\"],\n  [\"584\u00d75=2920\", \"243\u00d72=486\"],\n  [\"838\u00d75=4190\", \"711\u00d76=4266\"],\n  [\"106\u00d74=424\", \"412\u00d74=1648\"],\n  [\"634\u00d75=3170\", \"670\u00d74=2680\"],\n  [\"464\u00d75=2320\", \"200\u00d77=1400\"],\n  [\"648\u00d76=3888\", \"967\u00d73=2901\"],\n  [\"515\u00d73=1545\", \"921\u00d73=2763\"],\n  [\"465\u00d76=2790\", \"900\u00d75=4500\"],\n  [\"144\u00d78=1152\", \"550\u00d75=2750\"],\n  [\"904\u00d74=3616\", \"174\u00d76=1044\"],\n  [\"178\u00d77=1246\", \"332\u00d77=2324\"],\n  [\"427\u00d74=1708\", \"673\u00d74=2692\"],\n  [\"742\u00d77=5194\", \"403\u00d74=1612\"],\n  [\"488\u00d74=1952\", \"991\u00d76=5946\"],\n  [\"401\u00d78=3208\", \"663\u00d75=3315\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found in document: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Commit: \"Update master to output generated at 9a8706d\"\n# The document is a single table of \"A x B = C\" multiplication-fact cells.\n# This swaps 25 of those problem/answer strings for new ones while leaving\n# every other paragraph, row, and run formatting untouched.\n#\n# Each pair below is a unique, exact string that occurs exactly once in the\n# document, so a plain (non-wildcard) Find/Replace on each one is\n# unambiguous and preserves the run's existing formatting (rFonts/sz).\n\n$replacements = @(\n    @('651\u00d77=4557', '316\u00d74=1264'),\n    @('624\u00d77=4368', '190\u00d77=1330'),\n    @('144\u00d73=432', '152\u00d75=760'),\n    @('291\u00d78=2328', '129\u00d76=774'),\n    @('675\u00d76=4050', '212\u00d72=424'),\n    @('321\u00d76=1926', '959\u00d78=7672'),\n    @('200\u00d74=800', '514\u00d79=4626'),\n    @('389\u00d74=1556', '261\u00d74=1044'),\n    @('302\u00d78=2416', '521\u00d74=2084'),\n    @('873\u00d73=2619', '120\u00d73=360'),\n    @('584\u00d75=2920', '243\u00d72=486'),\n    @('838\u00d75=4190', '711\u00d76=4266'),\n    @('106\u00d74=424', '412\u00d74=1648'),\n    @('634\u00d75=3170', '670\u00d74=2680'),\n    @('464\u00d75=2320', '200\u00d77=1400'),\n    @('648\u00d76=3888', '967\u00d73=2901'),\n    @('515\u00d73=1545', '921\u00d73=2763'),\n    @('465\u00d76=2790', '900\u00d75=4500'),\n    @('144\u00d78=1152', '550\u00d75=2750'),\n    @('904\u00d74=3616', '174\u00d76=1044'),\n    @('178\u00d77=1246', '332\u00d77=2324'),\n    @('427\u00d74=1708', '673\u00d74=2692'),\n    @('742\u00d77=5194', '403\u00d74=1612'),\n    @('488\u00d74=1952', '991\u00d76=5946'),\n    @('401\u00d78=3208', '663\u00d75=3315')\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n\n    # wdFindContinue = 1 (Wrap), wdReplaceAll = 2 (Replace)\n    $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $found) {\n        throw \"Text not found in document: $oldText\"\n    }\n}\n"}
